$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column-B values (one per data row, r=2..47) representing the newly
# inserted most-recent observation. Existing B:K values shift right by one
# column and the previous column-K value is discarded (pushed out of range).
$newFirstCol = @{
    2 = -0.8260246319213993
    3 = -0.1405287498260583
    4 = 0.7815531178611421
    5 = 0.59085446968239
    6 = 0.9421746901075616
    7 = -0.8773451131212686
    8 = 0.04951982135029001
    9 = 0.2874422675287736
    10 = -0.5938852086106097
    11 = -0.2740769919852016
    12 = -0.09712588908289738
    13 = 0.1560540901775642
    14 = -0.106237102096024
    15 = -0.2673832081489602
    16 = -0.4749660705519546
    17 = 1.097015279354651
    18 = -0.2415579873788807
    19 = 1.205589761734299
    20 = -0.006071047505593896
    21 = -0.2580915896621678
    22 = 0.3657676764542774
    23 = 0.04413770072197692
    24 = 0.6424403654065582
    25 = 2.297389002388887
    26 = 8.826710628892494
    27 = -9.780318414391347
    28 = -1.200275438764269
    29 = 0.3719860057927588
    30 = -2.702915518772638
    31 = -0.2307826431404359
    32 = -0.565438627693374
    33 = -0.6603092772102132
    34 = -0.15162438770796
    35 = -0.2053460154962278
    36 = 0.6162032393936197
    37 = 1.652643173475852
    38 = 0.3110387314724781
    39 = 0.2388379152847414
    40 = 0.6508000635779043
    41 = 0.2387740594105157
    42 = 0.3465902496671606
    43 = 0.00230005330798793
    44 = -0.1902738424076751
    45 = -0.3325070745318338
    46 = 0.1656141382254278
    47 = -0.09587373626955231
}

for ($r = 2; $r -le 47; $r++) {
    # Capture the row's current B:K values before overwriting anything.
    # (Value2 is used for the read because it returns the raw numeric
    # value; Value can return a formatted/variant wrapper on read.)
    $old = @()
    for ($c = 2; $c -le 11; $c++) {
        $old += ,$ws.Cells.Item($r, $c).Value2
    }

    # Shift right: new C = old B, new D = old C, ... new K = old J.
    # (old K is dropped, matching the diff.)
    for ($c = 11; $c -ge 3; $c--) {
        $srcVal = $old[$c - 3]
        if ($null -eq $srcVal) {
            $ws.Cells.Item($r, $c).ClearContents()
        } else {
            $ws.Cells.Item($r, $c).Value = $srcVal
        }
    }

    # Insert the new value at column B.
    $ws.Cells.Item($r, 2).Value = $newFirstCol[$r]
}
